# Refresh the cryptos price/volume table (Coin/Link/Price/Volume(1h) columns)
# with the latest scraped figures. All these columns hold plain text in the
# source data (e.g. "27.247.16", "  +0.17%  "), so we force text formatting
# before the write and reset the style right after - this stops Excel's
# automatic "looks like a number" coercion (which would silently turn
# "1.000" into 1, or drop a display string's formatting) while still
# leaving the cell with no explicit style, matching the original cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddr, $NewValue)
    $r = $ws.Range($CellAddr)
    $r.NumberFormat = "@"
    $r.Value = $NewValue
    $r.Style = "Normal"
}

Set-TextValue 'D2' '27.247.16'
Set-TextValue 'E2' '  +0.17%  '
Set-TextValue 'D3' '1.907.76'
Set-TextValue 'E3' '  +0.13%  '
Set-TextValue 'E4' '  -0.12%  '
Set-TextValue 'D5' '307.52'
Set-TextValue 'E5' '  -0.11%  '
Set-TextValue 'E6' '  -0.08%  '
Set-TextValue 'D7' '0.5261'
Set-TextValue 'E7' '  +1.16%  '
Set-TextValue 'D8' '0.3812'
Set-TextValue 'E8' '  +1.19%  '
Set-TextValue 'D9' '0.07283'
Set-TextValue 'E9' '  +0.14%  '
Set-TextValue 'D10' '21.99'
Set-TextValue 'E10' '  +3.80%  '
Set-TextValue 'D11' '0.9023'
Set-TextValue 'E11' '  -0.26%  '
Set-TextValue 'D12' '0.08164'
Set-TextValue 'E12' '  -3.43%  '
Set-TextValue 'D13' '96.20'
Set-TextValue 'E13' '  -0.70%  '
Set-TextValue 'D14' '5.361'
Set-TextValue 'E14' '  +1.17%  '
Set-TextValue 'D15' '1.458.82'
Set-TextValue 'E15' '  -23.81%  '
Set-TextValue 'D16' '1.000'
Set-TextValue 'E16' '  -0.23%  '
Set-TextValue 'D17' '0.000008651'
Set-TextValue 'E17' '  -0.16%  '
Set-TextValue 'E18' '  +1.55%  '
Set-TextValue 'D20' '27.284.90'
Set-TextValue 'D21' '5.121'
Set-TextValue 'E21' '  +0.42%  '
Set-TextValue 'D23' '6.518'
Set-TextValue 'E23' '  +1.15%  '
Set-TextValue 'D24' '149.91'
Set-TextValue 'E24' '  +2.06%  '
Set-TextValue 'D25' '2.301'
Set-TextValue 'E25' '  -1.77%  '
Set-TextValue 'E26' '  +0.03%  '
Set-TextValue 'D27' '1.739'
Set-TextValue 'E27' '  -1.02%  '
Set-TextValue 'D28' '116.91'
Set-TextValue 'E28' '  +1.51%  '
Set-TextValue 'D29' '4.849'
Set-TextValue 'E29' '  +0.53%  '
Set-TextValue 'D30' '4.854'
Set-TextValue 'E30' '  -1.11%  '
Set-TextValue 'D31' '0.09246'
Set-TextValue 'E31' '  -0.39%  '
Set-TextValue 'D32' '0.8292'
Set-TextValue 'E32' '  +4.23%  '
Set-TextValue 'D33' '0.05070'
Set-TextValue 'E33' '  -0.15%  '
Set-TextValue 'D34' '1.227'
Set-TextValue 'E34' '  -1.27%  '
Set-TextValue 'D35' '2.989'
Set-TextValue 'E35' '  +1.41%  '
Set-TextValue 'D36' '3.349'
Set-TextValue 'E36' '  -2.24%  '
Set-TextValue 'D37' '2.714'
Set-TextValue 'E37' '  +5.01%  '
Set-TextValue 'D38' '0.5800'
Set-TextValue 'E38' '  -0.34%  '
Set-TextValue 'D39' '0.02006'
Set-TextValue 'E39' '  -0.08%  '
Set-TextValue 'D40' '1.078'
Set-TextValue 'E40' '  +0.18%  '
Set-TextValue 'D41' '9.232'
Set-TextValue 'E41' '  +1.78%  '
Set-TextValue 'D43' '116.75'
Set-TextValue 'E43' '  -0.02%  '
Set-TextValue 'D44' '0.1524'
Set-TextValue 'E44' '  +0.21%  '
Set-TextValue 'D45' '0.4930'
Set-TextValue 'E45' '  +0.69%  '
Set-TextValue 'B46' 'PaxDollar'
Set-TextValue 'C46' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D46' '1.000'
Set-TextValue 'E46' '  -0.08%  '
Set-TextValue 'B47' 'EnergySwap'
Set-TextValue 'C47' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D47' '10.16'
Set-TextValue 'E47' '  -0.08%  '
Set-TextValue 'D48' '1.647'
Set-TextValue 'E48' '  +0.51%  '
Set-TextValue 'D49' '39.02'
Set-TextValue 'E49' '  +3.41%  '
Set-TextValue 'D50' '0.06130'
Set-TextValue 'E50' '  +2.79%  '
Set-TextValue 'D51' '64.58'
Set-TextValue 'E51' '  +0.75%  '
